$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add "Sheet2" after the existing last sheet (matches the diff: new sheet
# appended at the end with sheetId="2" r:id="rId2").
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Fill in the summary/detail info for the process-scheduling run.
$ws2.Range("A1").Value = "Total Process"
$ws2.Range("B1").Value = 100

$ws2.Range("A2").Value = "AWT (Average Waiting Time)"
$ws2.Range("B2").Value = 487.81

$ws2.Range("A3").Value = "Total Waiting Time"
$ws2.Range("B3").Value = 48781

$ws2.Range("A4").Value = "ATAT (Average Turn Around Time)"
$ws2.Range("B4").Value = 500.39

$ws2.Range("A5").Value = "Total Turn Around Time"
$ws2.Range("B5").Value = 50039

# Restore Sheet1 as the active sheet (matches original tabSelected state).
$ws1.Activate()
